$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97-114 down to 98-115.
$ws.Rows.Item(97).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 97 with the new weekly record.
$ws.Cells.Item(97, 1).Value2  = 5
$ws.Cells.Item(97, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(97, 3).Value2  = "Maule"
$ws.Cells.Item(97, 4).Value2  = 44617
$ws.Cells.Item(97, 5).Value2  = 7
$ws.Cells.Item(97, 6).Value2  = "Fruta"
$ws.Cells.Item(97, 7).Value2  = 100108
$ws.Cells.Item(97, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(97, 9).Value2  = 100108002
$ws.Cells.Item(97, 10).Value2 = "Mango"
$ws.Cells.Item(97, 11).Value2 = "Sin especificar"
$ws.Cells.Item(97, 12).Value2 = "Primera"
$ws.Cells.Item(97, 13).Value2 = 324
$ws.Cells.Item(97, 14).Value2 = 6000
$ws.Cells.Item(97, 15).Value2 = 7000
$ws.Cells.Item(97, 16).Value2 = 6309
$ws.Cells.Item(97, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(97, 18).Value2 = "Perú"
$ws.Cells.Item(97, 19).Value2 = 1577
$ws.Cells.Item(97, 20).Value2 = 4
